$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.507.34"
$ws.Range("E2").Value = "  -3.23%  "
$ws.Range("D3").Value = "1.961.75"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "'322.12"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "'0.4771"
$ws.Range("E7").Value = "  -4.62%  "
$ws.Range("D8").Value = "'0.4075"
$ws.Range("E8").Value = "  -3.66%  "
$ws.Range("D9").Value = "'53.43"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "'0.08482"
$ws.Range("E10").Value = "  -4.88%  "
$ws.Range("E11").Value = "  -4.38%  "
$ws.Range("D12").Value = "'22.37"
$ws.Range("E12").Value = "  -3.88%  "
$ws.Range("D13").Value = "1.943.22"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").Value = "'7.646"
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("D15").Value = "'6.174"
$ws.Range("E15").Value = "  -4.31%  "
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "'90.14"
$ws.Range("E17").Value = "  -4.20%  "
$ws.Range("E18").Value = "  -3.43%  "
$ws.Range("D19").Value = "'0.06623"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").Value = "'18.63"
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "'5.833"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").Value = "28.519.91"
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("D25").Value = "'2.294"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "2.177.45"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("D27").Value = "'155.37"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "'5.943"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("D30").Value = "'2.167"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("D31").Value = "'124.00"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").Value = "'0.9873"
$ws.Range("E32").Value = "  -6.66%  "
$ws.Range("D33").Value = "'0.09604"
$ws.Range("E34").Value = "  -5.72%  "
$ws.Range("D35").Value = "'5.619"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").Value = "'3.664"
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("D37").Value = "'0.02355"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").Value = "'8.940"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "'0.06241"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").Value = "'0.6227"
$ws.Range("E41").Value = "  -4.76%  "
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  -5.44%  "
$ws.Range("D45").Value = "'1.358"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").Value = "'0.5981"
$ws.Range("E46").Value = "  -5.49%  "
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").Value = "'2.074"
$ws.Range("E48").Value = "  -5.95%  "
$ws.Range("D49").Value = "'3.404"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "'0.00000000329"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "'0.06838"
$ws.Range("E51").Value = "  -1.68%  "
